$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay text (column D is text-formatted in the source)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.245.35"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").Value = "3.506.59"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "584.17"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").Value = "135.08"
$ws.Range("E6").Value = "  +1.30%  "

$ws.Range("D7").Value = "3.506.56"
$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("D11").Value = "7.12"
$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("D12").Value = "0.377"
$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("D13").Value = "4.104.36"
$ws.Range("E13").Value = "  -0.41%  "

$ws.Range("D14").Value = "27.43"
$ws.Range("E14").Value = "  -0.80%  "

$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("E16").Value = "  +1.26%  "

$ws.Range("D17").Value = "3.507.14"
$ws.Range("E17").Value = "  -0.45%  "

$ws.Range("D18").Value = "64.264.88"
$ws.Range("E18").Value = "  -0.95%  "

$ws.Range("D19").Value = "9.78"
$ws.Range("E19").Value = "  -2.02%  "

$ws.Range("D20").Value = "13.88"
$ws.Range("E20").Value = "  -2.53%  "

$ws.Range("D21").Value = "5.59"
$ws.Range("E21").Value = "  -1.20%  "

$ws.Range("D22").Value = "383.96"
$ws.Range("E22").Value = "  -1.58%  "

$ws.Range("E23").Value = "  -1.31%  "

$ws.Range("D24").Value = "3.649.52"
$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("D25").Value = "73.88"
$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("E28").Value = "  +5.64%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").Value = "7.59"
$ws.Range("E30").Value = "  +2.05%  "

$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("D32").Value = "8.32"
$ws.Range("E32").Value = "  +1.06%  "

$ws.Range("E33").Value = "  -2.26%  "

$ws.Range("D34").Value = "3.519.74"
$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +0.49%  "

$ws.Range("D37").Value = "23.61"
$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("D38").Value = "5.32"
$ws.Range("E38").Value = "  +2.60%  "

$ws.Range("E39").Value = "  -2.91%  "

$ws.Range("D40").Value = "6.88"
$ws.Range("E40").Value = "  -0.90%  "

$ws.Range("D41").Value = "163.56"
$ws.Range("E41").Value = "  -4.84%  "

$ws.Range("D42").Value = "0.0782"
$ws.Range("E42").Value = "  -3.20%  "

$ws.Range("E43").Value = "  -0.92%  "

$ws.Range("D44").Value = "26.38"
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").Value = "1.22"
$ws.Range("E46").Value = "  -0.98%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "41.80"
$ws.Range("E47").Value = "  -0.94%  "

$ws.Range("D48").Value = "4.39"
$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("D49").Value = "1.60"
$ws.Range("E49").Value = "  -3.53%  "

$ws.Range("D50").Value = "2.478.50"
$ws.Range("E50").Value = "  -0.68%  "

$ws.Range("D51").Value = "6.76"
$ws.Range("E51").Value = "  -1.45%  "
